$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 'Prerequisites'
$ws.Range("D1").Value = 'Corequisites'
$ws.Range("E1").Value = 'Concurrent'
$ws.Range("F1").Value = 'Recommended'
$ws.Range("G1").Value = 'Terms Typically Offered'

$ws.Range("C2").Value = 'NA'
$ws.Range("D2").Value = 'NA'
$ws.Range("E2").Value = 'NA'
$ws.Range("F2").Value = 'NA'
$ws.Range("G2").Value = 'F, W, SP'

$ws.Range("C3").Value = 'NA'
$ws.Range("D3").Value = 'NA'
$ws.Range("E3").Value = 'NA'
$ws.Range("F3").Value = 'CM 102.'
$ws.Range("G3").Value = 'F, W, SP'

$ws.Range("C4").Value = 'NA'
$ws.Range("D4").Value = 'CM 113.'
$ws.Range("E4").Value = 'NA'
$ws.Range("F4").Value = 'CM 102.'
$ws.Range("G4").Value = 'F, W, SP '

$ws.Range("C5").Value = 'ARCE 106 or CM 113; MATH 141; and PHYS 141.'
$ws.Range("D5").Value = 'NA'
$ws.Range("E5").Value = 'NA'
$ws.Range("F5").Value = 'NA'
$ws.Range("G5").Value = 'F, W, SP'

$ws.Range("C6").Value = 'CM 115, PHYS 132 or CHEM 124.'
$ws.Range("D6").Value = 'CM 232.'
$ws.Range("E6").Value = 'NA'
$ws.Range("F6").Value = 'NA'
$ws.Range("G6").Value = 'F, W, SP '

$ws.Range("C7").Value = 'MATH 142 or MATH 182.'
$ws.Range("D7").Value = 'NA'
$ws.Range("E7").Value = 'NA'
$ws.Range("F7").Value = 'NA'
$ws.Range("G7").Value = 'F, W, SP'

$ws.Range("C8").Value = 'MATH 119 or equivalent.'
$ws.Range("D8").Value = 'NA'
$ws.Range("E8").Value = 'NA'
$ws.Range("F8").Value = 'NA'
$ws.Range("G8").Value = 'TBD'

$ws.Range("C9").Value = 'CE 113 or CM 115.'
$ws.Range("D9").Value = 'NA'
$ws.Range("E9").Value = 'NA'
$ws.Range("F9").Value = 'NA'
$ws.Range("G9").Value = 'F, W, SP'

$ws.Range("C10").Value = 'CM 113 or CE 259 or ARCE 106.'
$ws.Range("D10").Value = 'NA'
$ws.Range("E10").Value = 'NA'
$ws.Range("F10").Value = 'NA'
$ws.Range("G10").Value = 'F, SP'

$ws.Range("C11").Value = 'CM 214 and ARCE 212.'
$ws.Range("D11").Value = 'NA'
$ws.Range("E11").Value = 'NA'
$ws.Range("F11").Value = 'NA'
$ws.Range("G11").Value = 'F, W, SP'

$ws.Range("C12").Value = 'CM 239 or BRAE 239; CM 313.'
$ws.Range("D12").Value = 'CM 334.'
$ws.Range("E12").Value = 'NA'
$ws.Range("F12").Value = 'NA'
$ws.Range("G12").Value = 'F, W, SP '

$ws.Range("C13").Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of GE Area B1 with a grade of C- or better in at least one of the courses; and completion of GE Areas B2, B3, and B4.'
$ws.Range("D13").Value = 'NA'
$ws.Range("E13").Value = 'NA'
$ws.Range("F13").Value = 'NA'
$ws.Range("G13").Value = 'F,W,SP,SU'

$ws.Range("C14").Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'
$ws.Range("D14").Value = 'NA'
$ws.Range("E14").Value = 'NA'
$ws.Range("F14").Value = 'NA'
$ws.Range("G14").Value = 'F, W, SP'

$ws.Range("C15").Value = 'CM 115 and BUS 207.'
$ws.Range("D15").Value = 'NA'
$ws.Range("E15").Value = 'NA'
$ws.Range("F15").Value = 'NA'
$ws.Range("G15").Value = 'F, W, SP'

$ws.Range("C16").Value = 'BUS 215 and CM 232.'
$ws.Range("D16").Value = 'NA'
$ws.Range("E16").Value = 'NA'
$ws.Range("F16").Value = 'NA'
$ws.Range("G16").Value = 'F, W, SP'

$ws.Range("C17").Value = 'ARCE 106, CE 259 or CM 113.'
$ws.Range("D17").Value = 'NA'
$ws.Range("E17").Value = 'NA'
$ws.Range("F17").Value = 'NA'
$ws.Range("G17").Value = 'F, W, SP'

$ws.Range("C18").Value = 'Consent of instructor.'
$ws.Range("D18").Value = 'NA'
$ws.Range("E18").Value = 'NA'
$ws.Range("F18").Value = 'NA'
$ws.Range("G18").Value = 'TBD'

$ws.Range("C19").Value = 'CM 313.'
$ws.Range("D19").Value = 'NA'
$ws.Range("E19").Value = 'NA'
$ws.Range("F19").Value = 'NA'
$ws.Range("G19").Value = 'F, W, SP'

$ws.Range("C20").Value = 'CM 313.'
$ws.Range("D20").Value = 'NA'
$ws.Range("E20").Value = 'NA'
$ws.Range("F20").Value = 'NA'
$ws.Range("G20").Value = 'F,W,SP,SU'

$ws.Range("C21").Value = 'CM 413 and CM 480.'
$ws.Range("D21").Value = 'NA'
$ws.Range("E21").Value = 'NA'
$ws.Range("F21").Value = 'NA'
$ws.Range("G21").Value = 'TBD'

$ws.Range("C22").Value = 'Third-year standing.'
$ws.Range("D22").Value = 'NA'
$ws.Range("E22").Value = 'NA'
$ws.Range("F22").Value = 'NA'
$ws.Range("G22").Value = 'SP'

$ws.Range("C23").Value = 'Third-year standing.'
$ws.Range("D23").Value = 'NA'
$ws.Range("E23").Value = 'NA'
$ws.Range("F23").Value = 'NA'
$ws.Range("G23").Value = 'TBD'

$ws.Range("C24").Value = 'Third-year standing.'
$ws.Range("D24").Value = 'NA'
$ws.Range("E24").Value = 'NA'
$ws.Range("F24").Value = 'NA'
$ws.Range("G24").Value = 'F, W, SP'

$ws.Range("C25").Value = 'Third-year standing.'
$ws.Range("D25").Value = 'NA'
$ws.Range("E25").Value = 'NA'
$ws.Range("F25").Value = 'NA'
$ws.Range("G25").Value = 'TBD'

$ws.Range("C26").Value = 'Third-year standing.'
$ws.Range("D26").Value = 'NA'
$ws.Range("E26").Value = 'NA'
$ws.Range("F26").Value = 'NA'
$ws.Range("G26").Value = 'SP'

$ws.Range("C27").Value = 'Third-year standing.'
$ws.Range("D27").Value = 'NA'
$ws.Range("E27").Value = 'NA'
$ws.Range("F27").Value = 'NA'
$ws.Range("G27").Value = 'W'

$ws.Range("C28").Value = 'Third-year standing.'
$ws.Range("D28").Value = 'NA'
$ws.Range("E28").Value = 'NA'
$ws.Range("F28").Value = 'NA'
$ws.Range("G28").Value = 'TBD'

$ws.Range("C29").Value = 'CM 334; CM 335; and CM 413.'
$ws.Range("D29").Value = 'NA'
$ws.Range("E29").Value = 'NA'
$ws.Range("F29").Value = 'NA'
$ws.Range("G29").Value = 'F,W,SP,SU'

$ws.Range("C30").Value = 'STAT 251 or STAT 312, CM 313 and CM 334.'
$ws.Range("D30").Value = 'NA'
$ws.Range("E30").Value = 'NA'
$ws.Range("F30").Value = 'NA'
$ws.Range("G30").Value = 'F,W,SP,SU'

$ws.Range("C31").Value = 'CM 313; junior standing; Construction Management majors only.'
$ws.Range("D31").Value = 'NA'
$ws.Range("E31").Value = 'NA'
$ws.Range("F31").Value = 'NA'
$ws.Range("G31").Value = 'F, W, SP'

$ws.Range("C32").Value = 'CM 460 and consent of project advisor. See department for additional guidelines and requirements.'
$ws.Range("D32").Value = 'NA'
$ws.Range("E32").Value = 'NA'
$ws.Range("F32").Value = 'NA'
$ws.Range("G32").Value = 'F, SP'

$ws.Range("C33").Value = 'CM 460 and consent of project advisor. See department for additional guidelines and requirements.'
$ws.Range("D33").Value = 'NA'
$ws.Range("E33").Value = 'NA'
$ws.Range("F33").Value = 'NA'
$ws.Range("G33").Value = 'F, W, SP'

$ws.Range("C34").Value = 'CM 413.'
$ws.Range("D34").Value = 'CM 443.'
$ws.Range("E34").Value = 'NA'
$ws.Range("F34").Value = 'NA'
$ws.Range("G34").Value = 'TBD '

$ws.Range("C35").Value = 'Consent of instructor.'
$ws.Range("D35").Value = 'NA'
$ws.Range("E35").Value = 'NA'
$ws.Range("F35").Value = 'NA'
$ws.Range("G35").Value = 'TBD'

$ws.Range("C36").Value = 'Consent of instructor.'
$ws.Range("D36").Value = 'NA'
$ws.Range("E36").Value = 'NA'
$ws.Range("F36").Value = 'NA'
$ws.Range("G36").Value = 'TBD'

$ws.Range("C37").Value = 'Minimum junior standing.'
$ws.Range("D37").Value = 'NA'
$ws.Range("E37").Value = 'NA'
$ws.Range("F37").Value = 'NA'
$ws.Range("G37").Value = 'F, SP'

$ws.Range("C38").Value = 'CM 313.'
$ws.Range("D38").Value = 'NA'
$ws.Range("E38").Value = 'NA'
$ws.Range("F38").Value = 'NA'
$ws.Range("G38").Value = 'TBD'

$ws.Range("C39").Value = 'Consent of instructor.'
$ws.Range("D39").Value = 'NA'
$ws.Range("E39").Value = 'NA'
$ws.Range("F39").Value = 'NA'
$ws.Range("G39").Value = 'F, W, SP'

$ws.Range("C40").Value = 'Consent of instructor.'
$ws.Range("D40").Value = 'NA'
$ws.Range("E40").Value = 'NA'
$ws.Range("F40").Value = 'NA'
$ws.Range("G40").Value = 'TBD'

$ws.Range("C41").Value = 'Consent of instructor.'
$ws.Range("D41").Value = 'NA'
$ws.Range("E41").Value = 'NA'
$ws.Range("F41").Value = 'NA'
$ws.Range("G41").Value = 'TBD'

$ws.Range("C42").Value = 'Consent of program coordinator.'
$ws.Range("D42").Value = 'NA'
$ws.Range("E42").Value = 'NA'
$ws.Range("F42").Value = 'NA'
$ws.Range("G42").Value = 'TBD'

$ws.Range("C43").Value = 'Consent of program coordinator.'
$ws.Range("D43").Value = 'NA'
$ws.Range("E43").Value = 'NA'
$ws.Range("F43").Value = 'NA'
$ws.Range("G43").Value = 'TBD'

$ws.Range("C44").Value = 'Consent of program coordinator.'
$ws.Range("D44").Value = 'NA'
$ws.Range("E44").Value = 'NA'
$ws.Range("F44").Value = 'NA'
$ws.Range("G44").Value = 'TBD'

$ws.Range("C45").Value = 'Consent of program coordinator.'
$ws.Range("D45").Value = 'NA'
$ws.Range("E45").Value = 'NA'
$ws.Range("F45").Value = 'NA'
$ws.Range("G45").Value = 'TBD'

$ws.Range("C46").Value = 'Consent of program coordinator.'
$ws.Range("D46").Value = 'NA'
$ws.Range("E46").Value = 'NA'
$ws.Range("F46").Value = 'NA'
$ws.Range("G46").Value = 'TBD'

$ws.Range("C47").Value = 'Graduate standing or consent of instructor.'
$ws.Range("D47").Value = 'NA'
$ws.Range("E47").Value = 'NA'
$ws.Range("F47").Value = 'NA'
$ws.Range("G47").Value = 'TBD'

$ws.Range("C48").Value = 'Graduate standing or consent of instructor.'
$ws.Range("D48").Value = 'NA'
$ws.Range("E48").Value = 'NA'
$ws.Range("F48").Value = 'NA'
$ws.Range("G48").Value = 'TBD'

